# Add Week 15 simulations: update Rushing and Receiving stats tables.
$wb = $excel.ActiveWorkbook

# ---- Rushing sheet ----
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: C.Newton (index unchanged)
$rushing.Range("A2").Value = 0
$rushing.Range("C2").Value = 18
$rushing.Range("D2").Value = 16
$rushing.Range("E2").Value = 11
$rushing.Range("F2").Value = 7

# Row 3: P.Walker (index unchanged)
$rushing.Range("A3").Value = 1
$rushing.Range("C3").Value = 1
$rushing.Range("D3").Value = 3
$rushing.Range("E3").Value = 2
$rushing.Range("F3").Value = 0

# Row 4: C.Hubbard (index shifts 3 -> 2)
$rushing.Range("A4").Value = 2
$rushing.Range("C4").Value = 24
$rushing.Range("D4").Value = 15
$rushing.Range("E4").Value = 4
$rushing.Range("F4").Value = 8

# Row 5: R.Freeman (index shifts 4 -> 3, stats unchanged)
$rushing.Range("A5").Value = 3
$rushing.Range("C5").Value = 6
$rushing.Range("D5").Value = 3
$rushing.Range("E5").Value = 1
$rushing.Range("F5").Value = 0

# Row 6: A.Abdullah (index shifts 5 -> 4)
$rushing.Range("A6").Value = 4
$rushing.Range("C6").Value = 11
$rushing.Range("D6").Value = 10
$rushing.Range("E6").Value = 4
$rushing.Range("F6").Value = 3

# Row 7: Dj.Moore (index shifts 6 -> 5, stats unchanged)
$rushing.Range("A7").Value = 5
$rushing.Range("C7").Value = 1
$rushing.Range("D7").Value = 2
$rushing.Range("E7").Value = 1
$rushing.Range("F7").Value = 0

# Row 8: S.Smith (index shifts 7 -> 6, stats unchanged)
$rushing.Range("A8").Value = 6
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 1
$rushing.Range("E8").Value = 0
$rushing.Range("F8").Value = 0

# Row 9: T.Tremble (index shifts 8 -> 7, stats unchanged)
$rushing.Range("A9").Value = 7
$rushing.Range("C9").Value = 0
$rushing.Range("D9").Value = 1
$rushing.Range("E9").Value = 0
$rushing.Range("F9").Value = 1

# Row 10: I.Thomas (index shifts 9 -> 8, stats unchanged)
$rushing.Range("A10").Value = 8
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 0
$rushing.Range("E10").Value = 0
$rushing.Range("F10").Value = 0

# ---- Receiving sheet ----
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: C.Hubbard (index shifts 1 -> 0, stats unchanged)
$receiving.Range("A2").Value = 0
$receiving.Range("C2").Value = 15
$receiving.Range("D2").Value = 9
$receiving.Range("E2").Value = 2
$receiving.Range("F2").Value = 2
$receiving.Range("G2").Value = 0
$receiving.Range("H2").Value = 0

# Row 3: R.Freeman (index shifts 2 -> 1, stats unchanged)
$receiving.Range("A3").Value = 1
$receiving.Range("C3").Value = 6
$receiving.Range("D3").Value = 3
$receiving.Range("E3").Value = 0
$receiving.Range("F3").Value = 0
$receiving.Range("G3").Value = 0
$receiving.Range("H3").Value = 0

# Row 4: A.Abdullah (index shifts 3 -> 2)
$receiving.Range("A4").Value = 2
$receiving.Range("C4").Value = 23
$receiving.Range("D4").Value = 15
$receiving.Range("E4").Value = 1
$receiving.Range("F4").Value = 0
$receiving.Range("G4").Value = 2
$receiving.Range("H4").Value = 1

# Row 5: R.Anderson (index shifts 4 -> 3)
$receiving.Range("A5").Value = 3
$receiving.Range("C5").Value = 62
$receiving.Range("D5").Value = 32
$receiving.Range("E5").Value = 18
$receiving.Range("F5").Value = 3
$receiving.Range("G5").Value = 5
$receiving.Range("H5").Value = 2

# Row 6: Dj.Moore (index shifts 5 -> 4)
$receiving.Range("A6").Value = 4
$receiving.Range("C6").Value = 89
$receiving.Range("D6").Value = 58
$receiving.Range("E6").Value = 33
$receiving.Range("F6").Value = 14
$receiving.Range("G6").Value = 11
$receiving.Range("H6").Value = 6

# Row 7: T.Marshall (index shifts 6 -> 5, stats unchanged)
$receiving.Range("A7").Value = 5
$receiving.Range("C7").Value = 20
$receiving.Range("D7").Value = 13
$receiving.Range("E7").Value = 6
$receiving.Range("F7").Value = 1
$receiving.Range("G7").Value = 3
$receiving.Range("H7").Value = 2

# Row 8: B.Zylstra (index shifts 7 -> 6)
$receiving.Range("A8").Value = 6
$receiving.Range("C8").Value = 13
$receiving.Range("D8").Value = 12
$receiving.Range("E8").Value = 6
$receiving.Range("F8").Value = 4
$receiving.Range("G8").Value = 2
$receiving.Range("H8").Value = 2

# Row 9: S.Smith (index shifts 8 -> 7, stats unchanged)
$receiving.Range("A9").Value = 7
$receiving.Range("C9").Value = 4
$receiving.Range("D9").Value = 3
$receiving.Range("E9").Value = 1
$receiving.Range("F9").Value = 0
$receiving.Range("G9").Value = 1
$receiving.Range("H9").Value = 1

# Row 10: K.Kirkwood (index shifts 9 -> 8, stats unchanged)
$receiving.Range("A10").Value = 8
$receiving.Range("C10").Value = 4
$receiving.Range("D10").Value = 3
$receiving.Range("E10").Value = 2
$receiving.Range("F10").Value = 0
$receiving.Range("G10").Value = 1
$receiving.Range("H10").Value = 0

# Row 11: W.Snead (index shifts 10 -> 9, stats unchanged)
$receiving.Range("A11").Value = 9
$receiving.Range("C11").Value = 2
$receiving.Range("D11").Value = 1
$receiving.Range("E11").Value = 0
$receiving.Range("F11").Value = 0
$receiving.Range("G11").Value = 0
$receiving.Range("H11").Value = 0

# Row 12: A.Erickson (index shifts 11 -> 10, stats unchanged)
$receiving.Range("A12").Value = 10
$receiving.Range("C12").Value = 1
$receiving.Range("D12").Value = 1
$receiving.Range("E12").Value = 0
$receiving.Range("F12").Value = 0
$receiving.Range("G12").Value = 0
$receiving.Range("H12").Value = 0

# Row 13: T.Tremble (index shifts 12 -> 11)
$receiving.Range("A13").Value = 11
$receiving.Range("C13").Value = 19
$receiving.Range("D13").Value = 14
$receiving.Range("E13").Value = 5
$receiving.Range("F13").Value = 1
$receiving.Range("G13").Value = 3
$receiving.Range("H13").Value = 3

# Row 14: I.Thomas (index shifts 13 -> 12)
$receiving.Range("A14").Value = 12
$receiving.Range("C14").Value = 22
$receiving.Range("D14").Value = 11
$receiving.Range("E14").Value = 3
$receiving.Range("F14").Value = 2
$receiving.Range("G14").Value = 4
$receiving.Range("H14").Value = 2

# Row 15: C.Thompson (index shifts 14 -> 13, stats unchanged)
$receiving.Range("A15").Value = 13
$receiving.Range("C15").Value = 1
$receiving.Range("D15").Value = 0
$receiving.Range("E15").Value = 0
$receiving.Range("F15").Value = 0
$receiving.Range("G15").Value = 0
$receiving.Range("H15").Value = 0

Write-Host "Week 15 simulations added."
